# Bugs and broken links
#
# On the "4. Current content" slide, the sentence "Refer to Audit content
# for more detail." had its hyperlink only cover the words "Audit content",
# leaving the trailing text run as " for more detail." (leading space).
# The fix extends the hyperlink to also cover the single space that follows
# "Audit content", splitting that trailing run into a hyperlinked space run
# and a "for more detail." run (no longer starting with a space).

$p = $ppt.ActivePresentation

# Slides are 1-indexed and match the deck's visual order 1-13; slide 8 is
# "4. Current content" which contains the "Refer to Audit content for more
# detail." paragraph with the broken link boundary.
$slide = $p.Slides.Item(8)

# The paragraph lives in the single body placeholder shape on this slide.
$shape = $slide.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

# Locate the existing "Audit content" hyperlinked run so we can reuse its
# target address and find where it ends.
$linkRun = $tr.Find("Audit content", 0)
$address = $linkRun.ActionSettings(1).Hyperlink.Address

# The character immediately after "Audit content" is the space that should
# join the hyperlink run (currently it starts the " for more detail." run).
$spaceStart = $linkRun.Start + $linkRun.Length
$spaceRange = $tr.Characters($spaceStart, 1)

# Extending the hyperlink onto this single space splits the old
# " for more detail." run into a hyperlinked " " run and a "for more
# detail." run, matching the target markup.
$spaceRange.ActionSettings(1).Hyperlink.Address = $address
